$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 376968032.44
$ws.Range("P2").Value = 96961532.06
$ws.Range("Q2").Value = 47691798.33
$ws.Range("R2").Value = 1698.6688387562
$ws.Range("S2").Value = 74235893.83
$ws.Range("T2").Value = -11.7662379626
$ws.Range("U2").Value = 62101004.18
$ws.Range("V2").Value = 9.645814955000001
$ws.Range("W2").Value = 93985171.70999999
$ws.Range("X2").Value = 38642945.87
$ws.Range("Y2").Value = -10.1541289372
$ws.Range("Z2").Value = 704796.29
$ws.Range("AA2").Value = 100.059873512
$ws.Range("AB2").Value = 282982860.73
$ws.Range("AC2").Value = 35.375712644
$ws.Range("AD2").Value = 28.3286482575
$ws.Range("AE2").Value = 10.9403248076
$ws.Range("AF2").Value = 324.8255683068
$ws.Range("AG2").Value = 24.9318678567
